$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows for resource groups that should no longer be listed,
# keeping only rg-FinOps, rg-HPC and rg-test1 (rows 7, 8, 10 in the
# current sheet). Delete from the bottom up so row numbers above
# stay valid.
$ws.Rows.Item(9).Delete()   # rg-hub-spoke-lab-eastus
$ws.Rows.Item(6).Delete()   # rg-aspentech-eastus
$ws.Rows.Item(5).Delete()   # NetworkWatcherRG
$ws.Rows.Item(4).Delete()   # DefaultResourceGroup-EUS
$ws.Rows.Item(3).Delete()   # Default-ActivityLogAlerts
$ws.Rows.Item(2).Delete()   # cloud-shell-storage-eastus

# Reset the TagValue column to Value01, Value02, Value03.
$ws.Range("D2").Value = "Value01"
$ws.Range("D3").Value = "Value02"
$ws.Range("D4").Value = "Value03"

$ws.Range("D5").Select()
